$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing data (A:E) shifts to (B:F)
$ws.Columns.Item(1).Insert()

# New header for column B (old column A header position) -> "segments"
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)
$ws.Cells.Item(1, 2).Value = "segments"

# Fill new column A with 0-based segment index, matching style of column B (segment names, old col A)
for ($i = 0; $i -le 18; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $i
}

$excel.CutCopyMode = 0
